# Applies the "Updated cryptos list" price/volume refresh to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look numeric (e.g. "28.441.36", "0.4707") but are
# really text (coinranking often uses a "." thousands separator), so a leading
# apostrophe is used to force Excel to keep them as text instead of parsing them
# as numbers/dates, matching the original inlineStr cell content exactly.
$ws.Range("D2").Value = "'28.441.36"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "'1.871.44"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'339.44"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D7").Value = "'0.4707"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("D8").Value = "'0.3962"
$ws.Range("E8").Value = "  +3.61%  "
$ws.Range("D9").Value = "'47.34"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'0.08022"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").Value = "'1.003"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("D12").Value = "'21.91"
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("D13").Value = "'1.873.77"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "'5.996"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "'7.243"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").Value = "'91.42"
$ws.Range("E16").Value = "  +4.26%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").Value = "'0.06615"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'17.59"
$ws.Range("E20").Value = "  +4.05%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'28.416.59"
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("D23").Value = "'5.460"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").Value = "'11.05"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "'2.082.89"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").Value = "'160.37"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").Value = "'19.78"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "'2.136"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("D30").Value = "'5.522"
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").Value = "'120.12"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").Value = "'0.9692"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").Value = "'0.09492"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("D34").Value = "'3.568"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "'1.377"
$ws.Range("E35").Value = "  +4.49%  "
$ws.Range("D36").Value = "'5.358"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("D37").Value = "'0.06086"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("D38").Value = "'0.02248"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "'8.388"
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("D40").Value = "'1.184"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("D41").Value = "'0.5945"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "'0.1874"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("D44").Value = "'10.36"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("D46").Value = "'0.5586"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").Value = "'12.12"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").Value = "'1.957"
$ws.Range("E48").Value = "  +5.04%  "
$ws.Range("D49").Value = "'0.06856"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("D50").Value = "'2.079"
$ws.Range("E50").Value = "  +17.48%  "
$ws.Range("D51").Value = "'111.49"
$ws.Range("E51").Value = "  +1.56%  "
